# Update the dSF column (F) values to reflect repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 3
    3  = -4
    4  = 1
    5  = 0
    6  = -5
    7  = 3
    8  = 1
    9  = 2
    10 = -1
    11 = 1
    12 = -1
    13 = 6
    14 = -4
    15 = -7
    16 = 2
    18 = -1
    19 = 1
    20 = 7
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
